$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 656.6310342721592
$ws.Range("C3").Value = 554.3834728378716
$ws.Range("C4").Value = 554.2737198578341
$ws.Range("C5").Value = 626.0346221782509
$ws.Range("C6").Value = 624.95042417543
$ws.Range("C7").Value = 640.3493814850502
$ws.Range("C8").Value = 695.4750411678478
$ws.Range("C9").Value = 661.9225495763103
$ws.Range("C10").Value = 647.9990654226775
$ws.Range("C11").Value = 678.0768024284845
$ws.Range("C12").Value = 680.7148198624269
$ws.Range("C13").Value = 668.6534341868721
$ws.Range("C14").Value = 674.5582942518924
$ws.Range("C15").Value = 678.298759601792
$ws.Range("C16").Value = 696.7963625196249
$ws.Range("C17").Value = 709.9960680821971
$ws.Range("C18").Value = 727.7304293279441
$ws.Range("C19").Value = 722.5783429246624
$ws.Range("C20").Value = 724.0325490375233
$ws.Range("C21").Value = 729.9840552875692
$ws.Range("C22").Value = 739.2722755696913
$ws.Range("C23").Value = 749.4004912241836
$ws.Range("C24").Value = 759.1525413327415
$ws.Range("C25").Value = 765.0527908709126
$ws.Range("C26").Value = 771.3462583190944
$ws.Range("C27").Value = 777.2554163424372
$ws.Range("C28").Value = 786.0317567379705
$ws.Range("C29").Value = 791.2119179973587
$ws.Range("C30").Value = 794.5813229389697
$ws.Range("C31").Value = 801.0297373062248
$ws.Range("C32").Value = 804.7532427156723
$ws.Range("C33").Value = 808.3557613749459
$ws.Range("C34").Value = 811.5521788139013
$ws.Range("C35").Value = 814.0400577695273
$ws.Range("C36").Value = 816.2662377425578
$ws.Range("C37").Value = 818.9181048419148
$ws.Range("C38").Value = 820.3558913222331
$ws.Range("C39").Value = 822.4947884886146
$ws.Range("C40").Value = 825.6635991585881
$ws.Range("C41").Value = 828.1119404613956
$ws.Range("C42").Value = 828.4245631856591
$ws.Range("C43").Value = 830.2095209026251
$ws.Range("C44").Value = 831.4478019559867
$ws.Range("C45").Value = 832.5191975615543
$ws.Range("C46").Value = 832.8385203089094
$ws.Range("C47").Value = 833.6572875459192
$ws.Range("C48").Value = 834.914545806753
$ws.Range("C49").Value = 835.4956584563089
$ws.Range("C50").Value = 836.9161561460058
$ws.Range("C51").Value = 838.8917674043705
$ws.Range("C52").Value = 841.0826576311647
$ws.Range("C53").Value = 844.969708357901
$ws.Range("C54").Value = 848.6991714512504
$ws.Range("C55").Value = 854.5857068447127
$ws.Range("C56").Value = 860.2415293835763
$ws.Range("C57").Value = 862.7794304039149
$ws.Range("C58").Value = 863.8642467567737
$ws.Range("C59").Value = 865.1766039801112
$ws.Range("C60").Value = 869.6882321698185
$ws.Range("C61").Value = 879.0786329866407
$ws.Range("C62").Value = 880.3908458218101
$ws.Range("C63").Value = 881.6740978919713
$ws.Range("C64").Value = 882.637016271448
$ws.Range("C65").Value = 883.747866377127
